$d = $word.ActiveDocument

$replacements = @(
    @("77÷8=", "97÷2="),
    @("29÷3=", "54÷6="),
    @("51÷3=", "67÷3="),
    @("91÷2=", "96÷4="),
    @("75÷6=", "68÷5="),
    @("29÷6=", "21÷2="),
    @("24÷9=", "11÷3="),
    @("10÷2=", "28÷6="),
    @("85÷6=", "44÷9="),
    @("47÷3=", "14÷3="),
    @("18÷9=", "26÷6="),
    @("16÷5=", "25÷3="),
    @("67÷9=", "62÷5="),
    @("68÷6=", "49÷3="),
    @("12÷6=", "32÷3="),
    @("29÷7=", "78÷8="),
    @("93÷9=", "73÷6="),
    @("90÷3=", "36÷9="),
    @("79÷4=", "14÷7="),
    @("59÷5=", "57÷4="),
    @("52÷9=", "50÷3="),
    @("43÷7=", "98÷8="),
    @("12÷7=", "14÷2="),
    @("82÷2=", "97÷7="),
    @("13÷9=", "20÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
